# refactoring dca analysis into class structure
# Add a "Total Forecast" summary row (row 31) beneath the existing DCA rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Total Forecast"
$ws.Range("B31").Value = 672
$ws.Range("C31").Value = 268
$ws.Range("D31").Value = 6427
$ws.Range("E31").Value = 738
$ws.Range("F31").Value = 516

# Reflect the scrolled view (header row stays frozen, view scrolled down)
# and move the active selection to C35, matching the state captured after
# the edit.
$ws.Activate()
$ws.Range("C35").Select()
